# "changed team size in intro" -- update the "Team size:" line on the
# project-summary slide from " 4/5 members" to " ~4 members".

$p = $ppt.ActivePresentation

$targetOld = " 4/5 members"
$targetNew = " ~4 members"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shp = $slide.Shapes.Item($shi)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            $fullText = $tr.Text
            $idx = $fullText.IndexOf($targetOld)
            if ($idx -ge 0) {
                $start = $idx + 1
                $len = $targetOld.Length
                $sub = $tr.Characters($start, $len)
                $sub.Text = $targetNew
            }
        }
    }
}
